$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.281.12'
$ws.Range('E2').Value = '  -0.78%  '

$ws.Range('D3').Value = '1.622.04'
$ws.Range('E3').Value = '  -0.55%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.50%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.80'
$ws.Range('E6').Value = '  -1.34%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3789'
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.97'
$ws.Range('E8').Value = '  -2.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3606'
$ws.Range('E9').Value = '  -1.65%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.223'
$ws.Range('E10').Value = '  -5.07%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08060'
$ws.Range('E11').Value = '  -1.79%  '

$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.50%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.53'
$ws.Range('E13').Value = '  -3.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.528'
$ws.Range('E14').Value = '  -2.38%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001242'
$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.213'
$ws.Range('E16').Value = '  -3.44%  '

$ws.Range('D17').Value = '1.620.83'
$ws.Range('E17').Value = '  -0.25%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.34'
$ws.Range('E18').Value = '  -1.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06908'
$ws.Range('E19').Value = '  -0.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.86'
$ws.Range('E20').Value = '  -3.11%  '

$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.412'
$ws.Range('E22').Value = '  -2.89%  '

$ws.Range('D23').Value = '23.278.41'
$ws.Range('E23').Value = '  -0.82%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.71'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.181'
$ws.Range('E25').Value = '  +1.54%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.449'
$ws.Range('E26').Value = '  +0.57%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.03'
$ws.Range('E27').Value = '  -1.81%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.43'
$ws.Range('E28').Value = '  -1.25%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.282'
$ws.Range('E29').Value = '  -0.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.35'
$ws.Range('E30').Value = '  -1.64%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.292'
$ws.Range('E31').Value = '  -5.79%  '

$ws.Range('D32').Value = '1.802.53'
$ws.Range('E32').Value = '  -0.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.734'
$ws.Range('E33').Value = '  -2.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.88'
$ws.Range('E34').Value = '  +4.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9440'
$ws.Range('E35').Value = '  -3.76%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02783'
$ws.Range('E36').Value = '  -1.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2513'
$ws.Range('E37').Value = '  -0.97%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08816'
$ws.Range('E38').Value = '  -0.39%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.042'
$ws.Range('E39').Value = '  -3.74%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07093'
$ws.Range('E40').Value = '  -5.11%  '

$ws.Range('E41').Value = '  -3.37%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7004'
$ws.Range('E42').Value = '  -2.47%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.03'
$ws.Range('E43').Value = '  -1.64%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.23'
$ws.Range('E44').Value = '  -4.82%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6406'
$ws.Range('E46').Value = '  -3.43%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.305'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.983'
$ws.Range('E48').Value = '  -1.31%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07961'
$ws.Range('E49').Value = '  -0.84%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.194'
$ws.Range('E50').Value = '  -2.09%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.16'
$ws.Range('E51').Value = '  -5.71%  '
